# Generate Report for Handoff
#
# The localization run finished: rows that were showing the previous
# "Handed back: in sync with en-US" status are now "Ready for handoff",
# and the associated timestamps advance to the new report-generation time.
# Excel also re-shrinks the (now-narrower) status/date columns to fit the
# new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# Target stored column width (per the canonical OOXML) is ~17.216 "characters".
# This runtime's ColumnWidth setter snaps to an MDW=6 pixel grid (stored
# widths land on (pixels+5)/6), so feed the grid input that lands closest
# to the target stored width rather than the raw target itself.
$newColWidth = 16.3333333333333

# --- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = "2017-02-09 09:33:20"
$ws.Columns.Item(5).ColumnWidth = $newColWidth
$ws.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2017-02-09 09:32:58"
$ws.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2017-02-09 09:33:20"
$ws.Columns.Item(3).ColumnWidth = $newColWidth
